$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-28 12:54:46"

# Update the timestamp column (O) for all data rows (2 through 393)
for ($r = 2; $r -le 393; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Update the productAriaLabel column (M) for the rows whose text changed
# (the " - Online kein Bestand" suffix was removed for these products)
$ws.Cells.Item(144, 13).Value = "Pasquier Milchbrötchen 10St 3.50 Schweizer Franken"
$ws.Cells.Item(198, 13).Value = "Pasquier Pitch Schokolade 8 Stück 4.50 Schweizer Franken"
$ws.Cells.Item(251, 13).Value = "Pasquier Schokobrötchen 16 Stück 8.50 Schweizer Franken"
